# F1 2026 Prediction Model Technical Documentation: v2.0 -> v3.0 update
# Applies the Kaggle ETL pipeline narrative + refreshed numbers across the
# document body and the two results tables.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "MISSING: $old"
    }
}

# --- Title page ---------------------------------------------------------
Replace-Text "Technical Documentation v2.0" "Technical Documentation v3.0"

# --- 1. Executive Summary -----------------------------------------------
Replace-Text "The model predicts Max Verstappen (32.8%), Lando Norris (21.6%), and George Russell (14.4%) as the top three championship contenders. The GBM position predictor achieves a Leave-One-Out RMSE of 1.99, meaning predictions are accurate to within approximately 2 championship positions." "The model predicts Max Verstappen (34.8%), Lando Norris (23.0%), and George Russell (12.7%) as the top three championship contenders. The GBM position predictor, trained on 270 driver-seasons loaded from the Kaggle F1 World Championship dataset, achieves a Leave-One-Out RMSE of 1.72, meaning predictions are accurate to within approximately 2 championship positions."

# --- 2.1 Historical Data --------------------------------------------------
Replace-Text "The training dataset contains 159 driver-season records spanning the hybrid era (2014-2025). This period was selected because the 2014 turbo-hybrid regulation change is the closest historical analog to 2026's power unit revolution, making pre-2014 data less relevant for modeling regulation-change dynamics." "The training dataset contains 270 driver-season records spanning the hybrid era (2014-2025). For 2014-2024, all data is loaded at runtime from the Kaggle F1 World Championship dataset via an ETL pipeline (kaggle_loader.py) that extracts championship standings, race results, qualifying performance, top-10 rates, and DNF rates directly from CSV files. This replaced 300+ lines of manually curated data with a reproducible data pipeline, expanding training records from 159 to 270 (+70%) by including all grid drivers rather than a curated top-10 subset. The 2025 season data is hardcoded from official FIA standings and Fox Sports."

# --- 3.1 Gradient Boosting Position Predictor -----------------------------
Replace-Text "A scikit-learn GradientBoostingRegressor trained on 19 features to predict championship position. Hyperparameters: n_estimators=100, max_depth=3, learning_rate=0.1, subsample=0.8. The model is evaluated using Leave-One-Out cross-validation (LOO-CV), achieving an average RMSE of 1.99 across all 159 records." "A scikit-learn GradientBoostingRegressor trained on 19 features to predict championship position. Hyperparameters: n_estimators=100, max_depth=3, learning_rate=0.1, subsample=0.8. The model is evaluated using Leave-One-Out cross-validation (LOO-CV), achieving an average RMSE of 1.72 across all 270 records. The 14% RMSE improvement (from 1.99 with 159 records) is primarily due to backfield calibration: the expanded dataset includes 20-25 drivers per year instead of 10-16, giving the GBM anchor points for P15-P22 finishers."

# --- 3.2 Elo Rating System -------------------------------------------------
Replace-Text "A custom Elo system tracks driver skill evolution across seasons. The implementation uses a K-factor of 32 for rapid adaptation, with seasonal updates based on championship position. Starting ratings are initialized at 1500, and the system processes all 159 driver-seasons chronologically. Final 2025 Elo ratings range from 1375 (Hulkenberg) to 1747 (Verstappen)." "A custom Elo system tracks driver skill evolution across seasons. The implementation uses a K-factor of 32 for rapid adaptation, with seasonal updates based on championship position. Starting ratings are initialized at 1500, and the system processes all 270 driver-seasons chronologically. Final 2025 Elo ratings range from 1375 (Hulkenberg) to 1710 (Verstappen). The expanded dataset (from 159 to 270 records) caused Elo ratings to drop 30-40 points across the board, as more opponents in the system spread ratings more realistically."

# --- 6. Ensemble Weighting (expected points example) -----------------------
Replace-Text "Expected points (Avg Pts) come exclusively from the MC simulator's 10,000-season average, not from the ensemble. This means a driver can have lower expected points but higher ensemble win probability if the Bayesian model favors them (e.g., Norris: 288 pts but 21.6% win probability, vs. Russell: 297 pts but 14.4%)." "Expected points (Avg Pts) come exclusively from the MC simulator's 10,000-season average, not from the ensemble. This means a driver can have lower expected points but higher ensemble win probability if the Bayesian model favors them (e.g., Norris: 293 pts but 23.0% win probability, vs. Russell: 291 pts but 12.7%)."

# --- 7.1 GBM Position Predictor evaluation ----------------------------------
Replace-Text "Leave-One-Out Cross-Validation RMSE: 1.99. This means the model predicts championship finishing position within approximately 2 places on average. Year-by-year breakdown: 2020 RMSE=2.13 (14 drivers), 2021 RMSE=1.26 (13 drivers), 2022 RMSE=2.58 (13 drivers)." "Leave-One-Out Cross-Validation RMSE: 1.72. This means the model predicts championship finishing position within approximately 2 places on average, a 14% improvement over the previous curated dataset (RMSE 1.99 with 159 records)."

# --- 7.2 Feature Ablation Studies -------------------------------------------
Replace-Text "Qualifying performance features (avg_quali_pos, q3_rate, front_row_rate) were tested and rejected. Adding them increased LOO RMSE from 1.98 to 2.08 due to multicollinearity with existing features (win_rate, podium_rate). Top-10 rate and DNF rate were intentionally added only to the MC strength formula, not the GBM, to avoid this issue." "Qualifying performance features (avg_quali_pos, q3_rate, front_row_rate) were tested and rejected. Adding them increased LOO RMSE from 1.98 to 2.08 due to multicollinearity with existing features (win_rate, podium_rate). Top-10 rate and DNF rate were intentionally added only to the MC strength formula, not the GBM, to avoid this issue. Additionally, race-by-race finishing position variance (CV) and grid-to-finish delta were analyzed from the Kaggle dataset as potential features but excluded from the model to avoid overfitting (19 features on 270 records is already near the practical limit). These metrics are instead presented as exploratory insights in the dashboard's Insights tab."

# --- 9. Key Insights & Narratives -------------------------------------------
Replace-Text "Verstappen's Path to P1: Despite Red Bull's unproven Ford engine (maturity 80, 4th best), Verstappen's unmatched historical dominance (Elo 1747, 52% win rate, 96% top-10 rate) combined with peak age (28) and 0% DNF reliability makes him the model's top pick. His driver talent component (24.4) is nearly double the next best." "Verstappen's Path to P1: Despite Red Bull's unproven Ford engine (maturity 80, 4th best), Verstappen's unmatched historical dominance (Elo 1710, 52% win rate, 96% top-10 rate) combined with peak age (28) and 0% DNF reliability makes him the model's top pick. His driver talent component (24.4) is nearly double the next best."

Replace-Text "Hamilton's Age Penalty: The 7-time champion drives the 2nd-best car (Ferrari, expert rating 92) but his age-prime penalty (-5.7 at age 41) plus his poor 2025 form (0 podiums) significantly limits his title chances to 4.3%. His veteran bonus (experience 5.2 + reg-change 7.5 = 12.7) partially offsets the decline, but not enough to overcome 13 years past peak age." "Hamilton's Age Penalty: The 7-time champion drives the 2nd-best car (Ferrari, expert rating 92) but his age-prime penalty (-5.7 at age 41) plus his poor 2025 form (0 podiums) significantly limits his title chances to 4.1%. His veteran bonus (experience 5.2 + reg-change 7.5 = 12.7) partially offsets the decline, but not enough to overcome 13 years past peak age."

# --- 10. Code Architecture ---------------------------------------------------
Replace-Text "data/historical_data.py — 159 driver-season records (2014-2025), constructor results, 2026 grid definitions, pre-season testing data, bookmaker odds, top-10 rate and DNF rate data." "data/historical_data.py — Combines Kaggle-loaded 2014-2024 data with hardcoded 2025 season data and 2026 grid definitions, pre-season testing data, bookmaker odds."

Replace-Text "data/qualifying_data.py — Qualifying performance data (avg position, Q3 rate, front row rate) for all driver-seasons." "data/kaggle_loader.py — ETL pipeline that loads 2014-2024 driver results, constructor results, top-10/DNF rates, and qualifying stats from Kaggle F1 World Championship CSV files. Produces 247 driver-season records from 9 CSV source files."

Replace-Text "dashboard/f1-2026-predictor.html — Self-contained interactive React dashboard with team color coding and multi-tab navigation." "dashboard/f1-2026-predictor.html — Self-contained interactive React dashboard with team color coding and three tabs: Drivers (WDC predictions with model toggle), Constructors (WCC predictions), and Insights (finishing consistency CV and grid-to-finish delta analysis from Kaggle race-by-race data)."

# --- 12. Version History -----------------------------------------------------
Replace-Text "v2.0 (Current): Added F1 experience (actual career years with sqrt diminishing returns), regulation-change veteran bonus (2.5 pts per reg change survived), age-prime performance curve (peak at 25-31, decline after 35), top-10 finish rate consistency bonus, and DNF rate reliability penalty. Updated dashboard and documentation." "v2.0: Added F1 experience (actual career years with sqrt diminishing returns), regulation-change veteran bonus (2.5 pts per reg change survived), age-prime performance curve (peak at 25-31, decline after 35), top-10 finish rate consistency bonus, and DNF rate reliability penalty."

# Append the new v3.0 (Current) paragraph after the v2.0 paragraph, matching
# that paragraph's formatting (handled automatically by InsertParagraphAfter).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "v3.0 (Current): Replaced hardcoded 2014-2024 data with Kaggle F1 World Championship ETL pipeline (kaggle_loader.py). Training data expanded from 159 to 270 driver-seasons (+70%), LOO RMSE improved from 1.99 to 1.72 (-14%). Added Insights tab to dashboard with finishing consistency (CV) and grid-to-finish delta analysis from race-by-race Kaggle data. Removed qualifying_data.py (now loaded from Kaggle). Updated README with data pipeline documentation."

# --- 8.1 World Drivers' Championship table -----------------------------------
$driversTable = $d.Tables.Item(2)

# Row 2: Max Verstappen
$driversTable.Cell(2, 4).Range.Text = "37.8"
$driversTable.Cell(2, 5).Range.Text = "31.8"
$driversTable.Cell(2, 6).Range.Text = "34.8"
$driversTable.Cell(2, 7).Range.Text = "325"

# Row 3: Lando Norris
$driversTable.Cell(3, 4).Range.Text = "18.3"
$driversTable.Cell(3, 5).Range.Text = "27.8"
$driversTable.Cell(3, 6).Range.Text = "23.0"
$driversTable.Cell(3, 7).Range.Text = "293"

# Row 4: George Russell
$driversTable.Cell(4, 4).Range.Text = "17.0"
$driversTable.Cell(4, 5).Range.Text = "8.3"
$driversTable.Cell(4, 6).Range.Text = "12.7"
$driversTable.Cell(4, 7).Range.Text = "291"

# Row 5: Charles Leclerc
$driversTable.Cell(5, 4).Range.Text = "16.1"
$driversTable.Cell(5, 5).Range.Text = "3.7"
$driversTable.Cell(5, 6).Range.Text = "9.9"
$driversTable.Cell(5, 7).Range.Text = "287"

# Row 6: Oscar Piastri
$driversTable.Cell(6, 4).Range.Text = "3.2"
$driversTable.Cell(6, 5).Range.Text = "20.3"
$driversTable.Cell(6, 6).Range.Text = "11.8"
$driversTable.Cell(6, 7).Range.Text = "232"

# Row 7: Lewis Hamilton
$driversTable.Cell(7, 4).Range.Text = "7.3"
$driversTable.Cell(7, 5).Range.Text = "0.9"
$driversTable.Cell(7, 6).Range.Text = "4.1"
$driversTable.Cell(7, 7).Range.Text = "258"

# Row 8: Carlos Sainz
$driversTable.Cell(8, 5).Range.Text = "1.8"
$driversTable.Cell(8, 6).Range.Text = "0.9"
$driversTable.Cell(8, 7).Range.Text = "137"

# Row 9: Kimi Antonelli
$driversTable.Cell(9, 4).Range.Text = "0.4"
$driversTable.Cell(9, 5).Range.Text = "1.1"
$driversTable.Cell(9, 7).Range.Text = "189"

# --- 8.2 World Constructors' Championship table -------------------------------
$constructorsTable = $d.Tables.Item(3)

$constructorsTable.Cell(2, 4).Range.Text = "545"   # Ferrari
$constructorsTable.Cell(3, 4).Range.Text = "525"   # McLaren
$constructorsTable.Cell(4, 4).Range.Text = "480"   # Mercedes
$constructorsTable.Cell(5, 4).Range.Text = "401"   # Red Bull
$constructorsTable.Cell(6, 4).Range.Text = "202"   # Williams

Write-Output "done"
